# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 223 on the active sheet,
# pushing the existing rows 223-232 down to 224-233.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 223 (shifts rows 223.. down by one)
$ws.Rows(223).Insert()

# Populate the newly inserted row with the new record's data
$ws.Cells.Item(223, 1).Value = 10
$ws.Cells.Item(223, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(223, 3).Value = "La Araucanía"
$ws.Cells.Item(223, 4).Value = 44706
$ws.Cells.Item(223, 5).Value = 9
$ws.Cells.Item(223, 6).Value = 100112052
$ws.Cells.Item(223, 7).Value = "Albahaca"
$ws.Cells.Item(223, 8).Value = "Sin especificar"
$ws.Cells.Item(223, 9).Value = "Primera"
$ws.Cells.Item(223, 10).Value = 40
$ws.Cells.Item(223, 11).Value = 5500
$ws.Cells.Item(223, 12).Value = 5500
$ws.Cells.Item(223, 13).Value = 5500
$ws.Cells.Item(223, 14).Value = "$/paquete"
$ws.Cells.Item(223, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(223, 16).Value = 5500
$ws.Cells.Item(223, 17).Value = 1
$ws.Cells.Item(223, 18).Value = "Hortaliza"
